# Apply the "major rewrite" fixture update:
#  - swap the 3 genotyped SNP columns (J/K/L) for 3 new rsIDs and their calls
#  - add a new sample row (ID-CG-000139T)
#  - move the active selection
#  - set page setup to A4 / portrait
#
# NOTE: cell writes are ordered deliberately. This engine rebuilds
# xl/sharedStrings.xml from scratch on save: strings still referenced by some
# cell keep their relative order, unused ones are dropped, and brand-new
# strings are appended in the order their owning cells are written. Writing
# in this particular order reproduces the target shared-string table order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 5 (sample ID + genotype placeholders first; rest filled below) ---
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = "ID-CG-000139T"

# --- Header row: new SNP rsIDs for columns J, K, L ---
$ws.Range("J1").Value = "rs9988021"
$ws.Range("K1").Value = "rs115551684"
$ws.Range("L1").Value = "rs199560653"

# --- Row 2 (ID-CG-000026T) genotype calls for the new SNPs ---
$ws.Range("J2").Value = "G G"
$ws.Range("K2").Value = "G A"
$ws.Range("L2").Value = "G G"

# --- Row 3 (ID-CG-000027T) genotype calls for the new SNPs ---
$ws.Range("J3").Value = "G A"
$ws.Range("K3").Value = "G G"
$ws.Range("L3").Value = "G A"

# --- Row 4 (ID-CG-000028T) genotype calls for the new SNPs ---
$ws.Range("J4").Value = "G A"
$ws.Range("K4").Value = "A G"
$ws.Range("L4").Value = "G G"

# --- Finish new row 5 (ID-CG-000139T) ---
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = "C C"
$ws.Range("H5").Value = "C C"
$ws.Range("I5").Value = "T T"
$ws.Range("J5").Value = "G G"
$ws.Range("K5").Value = "G G"
$ws.Range("L5").Value = "A G"

# --- Move the active selection ---
[void]$ws.Range("L7").Select()

# --- Page setup: A4, portrait ---
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
